$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column E: "along" header, with "cut_aid_in_programTRUE" for every data row
$ws.Range("E1").Value = "along"
for ($r = 2; $r -le 12; $r++) {
    $ws.Range("E$r").Value = "cut_aid_in_programTRUE"
}

# Updated numeric values for columns B, C, D (rows 2-12)
$values = @{
    "B2" = -0.0513286442643204;  "C2" = -0.0671790184472558;  "D2" = -0.0354782700813851
    "B3" = -0.0478584132505946;  "C3" = -0.0706141477727704;  "D3" = -0.0251026787284189
    "B4" = 0.0059073444719303;   "C4" = -0.0503131381821435;  "D4" = 0.0621278271260041
    "B5" = -0.0495346534007151;  "C5" = -0.0974637307836841;  "D5" = -0.00160557601774618
    "B6" = -0.124368455268675;   "C6" = -0.179542730281175;   "D6" = -0.0691941802561751
    "B7" = -0.072172816617106;   "C7" = -0.142611302671681;   "D7" = -0.0017343305625309
    "B8" = -0.0915766986174901;  "C8" = -0.154893116051925;   "D8" = -0.0282602811830557
    "B9" = -0.00537949765378601; "C9" = -0.0604664383018167;  "D9" = 0.0497074429942447
    "B10" = 0.0868622068295776;  "C10" = 0.0084108846117876;  "D10" = 0.165313529047368
    "B11" = -0.0148317041662747; "C11" = -0.0493904940584233; "D11" = 0.0197270857258738
    "B12" = -0.0692743648953208; "C12" = -0.096662491837547;  "D12" = -0.0418862379530946
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
